$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 105 (shifts existing rows 105-116 down to 106-117)
$ws.Rows(105).Insert()

# Populate the newly inserted row 105 with the new weekly record
$ws.Range("A105").Value = 10
$ws.Range("B105").Value = "Vega Modelo de Temuco"
$ws.Range("C105").Value = "La Araucanía"
$ws.Range("D105").Value = 45142
$ws.Range("E105").Value = 9
$ws.Range("F105").Value = "Fruta"
$ws.Range("G105").Value = 100108
$ws.Range("H105").Value = "Tropicales y subtropicales"
$ws.Range("I105").Value = 100108004
$ws.Range("J105").Value = "Papaya"
$ws.Range("K105").Value = "Cultivar IV Región"
$ws.Range("L105").Value = "Primera"
$ws.Range("M105").Value = 55
$ws.Range("N105").Value = 25000
$ws.Range("O105").Value = 25000
$ws.Range("P105").Value = 25000
$ws.Range("Q105").Value = "$/bandeja 10 kilos"
$ws.Range("R105").Value = "Provincia del Elquí"
$ws.Range("S105").Value = 2500
$ws.Range("T105").Value = 10
